$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 20833596
$ws.Range("I6").Value = 55555650
$ws.Range("K6").Value = 166666950
$ws.Range("M6").Value = -166666838
$ws.Range("H38").Value = 5962.8
$ws.Range("I38").Value = 1604.6666
$ws.Range("K38").Value = 4813.9998
$ws.Range("M38").Value = -4441.9998
$ws.Range("H40").Value = 3786.1428
$ws.Range("I40").Value = 2167.6667
$ws.Range("K40").Value = 2167.6667
$ws.Range("M40").Value = -1992.6667
$ws.Range("H51").Value = 35689.875
$ws.Range("I51").Value = 5232.5
$ws.Range("J51").Value = 40040.93
$ws.Range("K51").Value = 5232.5
$ws.Range("L51").Value = 40040.93
$ws.Range("M51").Value = -4748.5
$ws.Range("N51").Value = -41008.93
$ws.Range("H58").Value = 2559.7222
$ws.Range("I58").Value = 107.6
$ws.Range("K58").Value = 322.8
$ws.Range("M58").Value = -172.8
$ws.Range("H86").Value = 30571.625
$ws.Range("I86").Value = 2915.3
$ws.Range("J86").Value = 76665.5
$ws.Range("K86").Value = 2915.3
$ws.Range("L86").Value = 76665.5
$ws.Range("M86").Value = -1792.3
$ws.Range("N86").Value = -78911.5
$ws.Range("H89").Value = 30571.625
$ws.Range("I89").Value = 2915.3
$ws.Range("J89").Value = 76665.5
$ws.Range("K89").Value = 14576.5
$ws.Range("L89").Value = 383327.5
$ws.Range("M89").Value = -8960.5
$ws.Range("N89").Value = -394559.5
$ws.Range("H106").Value = 5832.9443
$ws.Range("I106").Value = 4857
$ws.Range("J106").Value = 9248.75
$ws.Range("K106").Value = 4857
$ws.Range("L106").Value = 9248.75
$ws.Range("M106").Value = -4226
$ws.Range("N106").Value = -10510.75
$ws.Range("H112").Value = 36908.266
$ws.Range("J112").Value = 40325.355
$ws.Range("L112").Value = 120976.065
$ws.Range("N112").Value = -123192.065
$ws.Range("H113").Value = 13329.947
$ws.Range("I113").Value = 20028.8
$ws.Range("J113").Value = 5886.778
$ws.Range("K113").Value = 20028.8
$ws.Range("L113").Value = 5886.778
$ws.Range("M113").Value = -16774.8
$ws.Range("N113").Value = -12394.778
$ws.Range("H116").Value = 4199.4287
$ws.Range("I116").Value = 3479.4
$ws.Range("J116").Value = 5999.5
$ws.Range("K116").Value = 3479.4
$ws.Range("L116").Value = 5999.5
$ws.Range("M116").Value = -37.40000000000009
$ws.Range("N116").Value = -12883.5
$ws.Range("H132").Value = 7664.467
$ws.Range("I132").Value = 11631.375
$ws.Range("J132").Value = 3130.8572
$ws.Range("K132").Value = 34894.125
$ws.Range("L132").Value = 9392.571599999999
$ws.Range("M132").Value = -32364.125
$ws.Range("N132").Value = -14452.5716
$ws.Range("H138").Value = 9720.034
$ws.Range("I138").Value = 10433
$ws.Range("J138").Value = 9344.789000000001
$ws.Range("K138").Value = 31299
$ws.Range("L138").Value = 28034.367
$ws.Range("M138").Value = -26159
$ws.Range("N138").Value = -38314.367
$ws.Range("H141").Value = 5547.6333
$ws.Range("I141").Value = 4497.381
$ws.Range("K141").Value = 13492.143
$ws.Range("M141").Value = -8312.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1206.2727
$ws.Range("I4").Value = 908.75
$ws.Range("J4").Value = 1999.6666
$ws.Range("K4").Value = 908.75
$ws.Range("L4").Value = 1999.6666
$ws.Range("M4").Value = -792.75
$ws.Range("N4").Value = -2231.6666
$ws.Range("H5").Value = 929.8461
$ws.Range("I5").Value = 228.875
$ws.Range("J5").Value = 2051.4
$ws.Range("K5").Value = 228.875
$ws.Range("L5").Value = 2051.4
$ws.Range("M5").Value = -116.875
$ws.Range("N5").Value = -2275.4
$ws.Range("H32").Value = 3282.6924
$ws.Range("I32").Value = 3352.2432
$ws.Range("K32").Value = 3352.2432
$ws.Range("M32").Value = -3065.2432
$ws.Range("H122").Value = 362298.3
$ws.Range("I122").Value = 2460.0715
$ws.Range("K122").Value = 7380.2145
$ws.Range("M122").Value = -4930.2145
$ws.Range("H132").Value = 4226.5
$ws.Range("I132").Value = 3737
$ws.Range("J132").Value = 5499.2
$ws.Range("K132").Value = 11211
$ws.Range("L132").Value = 16497.6
$ws.Range("M132").Value = -8681
$ws.Range("N132").Value = -21557.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 929.8461
$ws.Range("I4").Value = 228.875
$ws.Range("J4").Value = 2051.4
$ws.Range("K4").Value = 228.875
$ws.Range("L4").Value = 2051.4
$ws.Range("M4").Value = -113.875
$ws.Range("N4").Value = -2281.4
$ws.Range("H132").Value = 79333
$ws.Range("J132").Value = 79333
$ws.Range("L132").Value = 79333
$ws.Range("N132").Value = -89453
$ws.Range("H134").Value = 9049.579
$ws.Range("I134").Value = 9049.579
$ws.Range("K134").Value = 27148.737
$ws.Range("M134").Value = -24613.737

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2491.4792
$ws.Range("I31").Value = 1166.5834
$ws.Range("J31").Value = 2933.111
$ws.Range("K31").Value = 1166.5834
$ws.Range("L31").Value = 2933.111
$ws.Range("M31").Value = -871.5834
$ws.Range("N31").Value = -3523.111
$ws.Range("H34").Value = 2491.4792
$ws.Range("I34").Value = 1166.5834
$ws.Range("J34").Value = 2933.111
$ws.Range("K34").Value = 1166.5834
$ws.Range("L34").Value = 2933.111
$ws.Range("M34").Value = -964.5834
$ws.Range("N34").Value = -3337.111
$ws.Range("H58").Value = 5521.024
$ws.Range("I58").Value = 6526.16
$ws.Range("K58").Value = 6526.16
$ws.Range("M58").Value = -6323.16
$ws.Range("H107").Value = 16691.643
$ws.Range("I107").Value = 24720.445
$ws.Range("J107").Value = 2239.8
$ws.Range("K107").Value = 24720.445
$ws.Range("L107").Value = 2239.8
$ws.Range("M107").Value = -22800.445
$ws.Range("N107").Value = -6079.8
$ws.Range("H122").Value = 2749.5
$ws.Range("I122").Value = 2749.5
$ws.Range("K122").Value = 8248.5
$ws.Range("M122").Value = -5798.5
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 66771984
$ws.Range("I132").Value = 111151640
$ws.Range("J132").Value = 202500
$ws.Range("K132").Value = 333454920
$ws.Range("L132").Value = 607500
$ws.Range("M132").Value = -333452390
$ws.Range("N132").Value = -612560
$ws.Range("H136").Value = 5521.024
$ws.Range("I136").Value = 6526.16
$ws.Range("K136").Value = 19578.48
$ws.Range("M136").Value = -17028.48

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 651.1875
$ws.Range("I92").Value = 429.63635
$ws.Range("J92").Value = 1138.6
$ws.Range("K92").Value = 1288.90905
$ws.Range("L92").Value = 3415.8
$ws.Range("M92").Value = -40.90904999999998
$ws.Range("N92").Value = -5911.799999999999
$ws.Range("H132").Value = 3799697
$ws.Range("J132").Value = 4518486.5
$ws.Range("L132").Value = 40666378.5
$ws.Range("N132").Value = -40671438.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 11810.75
$ws.Range("I132").Value = 9415.333000000001
$ws.Range("K132").Value = 28245.999
$ws.Range("M132").Value = -25715.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22582
$ws.Range("I7").Value = 35998.46
$ws.Range("J7").Value = 5140.6
$ws.Range("K7").Value = 35998.46
$ws.Range("L7").Value = 5140.6
$ws.Range("M7").Value = -35886.46
$ws.Range("N7").Value = -5364.6
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H45").Value = 30041
$ws.Range("I45").Value = 30041
$ws.Range("K45").Value = 30041
$ws.Range("M45").Value = -29634
$ws.Range("H48").Value = 36759
$ws.Range("I48").Value = 14020.5
$ws.Range("K48").Value = 14020.5
$ws.Range("M48").Value = -13359.5
$ws.Range("H122").Value = 5659.375
$ws.Range("I122").Value = 4843.75
$ws.Range("J122").Value = 6475
$ws.Range("K122").Value = 14531.25
$ws.Range("L122").Value = 19425
$ws.Range("M122").Value = -12081.25
$ws.Range("N122").Value = -24325
$ws.Range("H125").Value = 66665.664
$ws.Range("J125").Value = 66665.664
$ws.Range("L125").Value = 66665.664
$ws.Range("N125").Value = -76505.664
$ws.Range("H126").Value = 22582
$ws.Range("I126").Value = 35998.46
$ws.Range("J126").Value = 5140.6
$ws.Range("K126").Value = 107995.38
$ws.Range("L126").Value = 15421.8
$ws.Range("M126").Value = -105525.38
$ws.Range("N126").Value = -20361.8
$ws.Range("H132").Value = 19272.438
$ws.Range("I132").Value = 31463.223
$ws.Range("J132").Value = 3598.5715
$ws.Range("K132").Value = 94389.66900000001
$ws.Range("L132").Value = 10795.7145
$ws.Range("M132").Value = -91859.66900000001
$ws.Range("N132").Value = -15855.7145
$ws.Range("H136").Value = 5359.303
$ws.Range("I136").Value = 1622.25
$ws.Range("J136").Value = 11108.615
$ws.Range("K136").Value = 4866.75
$ws.Range("L136").Value = 33325.845
$ws.Range("M136").Value = -2316.75
$ws.Range("N136").Value = -38425.845

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 51090.332
$ws.Range("J74").Value = 51090.332
$ws.Range("L74").Value = 51090.332
$ws.Range("N74").Value = -52962.332
$ws.Range("H77").Value = 51090.332
$ws.Range("J77").Value = 51090.332
$ws.Range("L77").Value = 153270.996
$ws.Range("N77").Value = -162630.996
$ws.Range("H132").Value = 38831.332
$ws.Range("I132").Value = 56397.816
$ws.Range("K132").Value = 169193.448
$ws.Range("M132").Value = -166663.448
